# Reorder "Recorded By" (column G) entries so that "System" moves from the
# front of the comma-separated list to the end, for every row where the
# cell value begins with "System, ".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "System, *") {
        $rest = $val.Substring(8)  # strip leading "System, "
        $newVal = $rest + ", System"
        $cell.Value = $newVal
    }
}
